$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to remain text
# (values like "0.9957" or "30.919.43" would otherwise be auto-coerced
# into numbers by Excel's type inference on Range.Value assignment).
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.919.43'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.947.39'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").Value = '0.9957'
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").Value = '244.89'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").Value = '0.9962'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '0.4870'
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("D8").Value = '0.2962'
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("D9").Value = '0.06824'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("D11").Value = '107.05'
$ws.Range("E11").Value = '  -4.49%  '
$ws.Range("D12").Value = '1.951.87'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '5.462'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = '0.7058'
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("D16").Value = '281.57'
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").Value = '30.941.43'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").Value = '13.24'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '0.000007736'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = '2.201.89'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = '0.9959'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").Value = '5.499'
$ws.Range("E22").Value = '  -2.64%  '
$ws.Range("D23").Value = '0.9949'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '6.493'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Value = '9.815'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").Value = '169.11'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '19.96'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").Value = '2.215'
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("E29").Value = '  -3.13%  '
$ws.Range("D30").Value = '1.409'
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("D31").Value = '1.582'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = '4.558'
$ws.Range("E32").Value = '  -2.24%  '
$ws.Range("D33").Value = '4.477'
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").Value = '0.04954'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("D35").Value = '0.7643'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").Value = '0.02026'
$ws.Range("E38").Value = '  -2.36%  '
$ws.Range("D39").Value = '2.689'
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("D40").Value = '2.154'
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("D41").Value = '6.503'
$ws.Range("E41").Value = '  +8.72%  '
$ws.Range("D42").Value = '75.24'
$ws.Range("E42").Value = '  +8.19%  '
$ws.Range("D43").Value = '0.4491'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").Value = '109.22'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("D45").Value = '0.8819'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("D46").Value = '8.158'
$ws.Range("E46").Value = '  +10.37%  '
$ws.Range("D47").Value = '0.9958'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '979.21'
$ws.Range("E48").Value = '  +7.73%  '
$ws.Range("D49").Value = '9.375'
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").Value = '0.1263'
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("D51").Value = '35.72'
$ws.Range("E51").Value = '  +0.10%  '
